$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "wong3"

$ws.Range("B2").Value = 28
$ws.Range("C2").Value = 25
$ws.Range("E2").Value = 84
$ws.Range("F2").Value = 83
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 393

$ws.Range("B4").Value = 24
$ws.Range("C4").Value = 23
$ws.Range("E4").Value = 11
$ws.Range("F4").Value = 10
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 393

$ws.Range("B8").Value = 264
$ws.Range("C8").Value = 264
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 55

$ws.Range("B12").Value = 196
$ws.Range("C12").Value = 195
$ws.Range("E12").Value = 92
$ws.Range("F12").Value = 91

$ws.Range("B13").Value = 19
$ws.Range("C13").Value = 19
$ws.Range("E13").Value = 62
$ws.Range("F13").Value = 62
$ws.Range("J13").Value = 12
$ws.Range("K13").Value = 409
